# This script consolidates the 384-well "Dest" plate numbering used by the
# Worklist and Platemap sheets: plates were numbered 384-1..384-8 (one plate
# per 16 source wells); the edit merges each consecutive pair of plates
# (384-1+384-2 -> 384-1, 384-3+384-4 -> 384-2, 384-5+384-6 -> 384-3,
# 384-7+384-8 -> 384-4). For the second plate of each pair (the even-numbered
# source plate), the destination well also moves into the "second half" of
# the merged plate: the numeric DestWell shifts by +176, and on the Platemap
# sheet the well-letter/column label's column number shifts from 03 to 14.

$wb = $excel.ActiveWorkbook

function Get-DestPlateNum($n) {
    return [int][Math]::Ceiling($n / 2.0)
}

# ---- Sheet "Worklist" (columns: A Source, B Well, C Dest, D DestWell, E Volume) ----
$ws1 = $wb.Worksheets.Item("Worklist")

for ($row = 18; $row -le 116; $row++) {
    $destCell = $ws1.Cells.Item($row, 3)   # column C
    $wellCell = $ws1.Cells.Item($row, 4)   # column D

    $destVal = [string]$destCell.Value2
    if ($destVal -match '^384-(\d+)$') {
        $n = [int]$matches[1]
        $newDest = "384-" + (Get-DestPlateNum $n)
        $destCell.Value = $newDest

        if ($n % 2 -eq 0) {
            $wellVal = [string]$wellCell.Value2
            $newWell = [string]([int]$wellVal + 176)
            $wellCell.Value = "'" + $newWell
        }
    }
}

# ---- Sheet "Platemap" (columns: A.. , D Dest plate, E Well#, F Well label) ----
$ws2 = $wb.Worksheets.Item("Platemap")

for ($row = 17; $row -le 115; $row++) {
    $destCell  = $ws2.Cells.Item($row, 4)  # column D
    $wellCell  = $ws2.Cells.Item($row, 5)  # column E
    $labelCell = $ws2.Cells.Item($row, 6)  # column F

    $destVal = [string]$destCell.Value2
    if ($destVal -match '^384-(\d+)$') {
        $n = [int]$matches[1]
        $newDest = "384-" + (Get-DestPlateNum $n)
        $destCell.Value = $newDest

        if ($n % 2 -eq 0) {
            $wellVal = [string]$wellCell.Value2
            $newWell = [string]([int]$wellVal + 176)
            $wellCell.Value = "'" + $newWell

            $labelVal = [string]$labelCell.Value2
            $letter = $labelVal.Substring(0, 1)
            $labelCell.Value = $letter + "14"
        }
    }
}
